$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 882.35297
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -886
$ws.Range("N4").Value = -728
$ws.Range("H21").Value = 19870.467
$ws.Range("I21").Value = 20773.615
$ws.Range("J21").Value = 14000
$ws.Range("K21").Value = 20773.615
$ws.Range("L21").Value = 14000
$ws.Range("M21").Value = -20305.615
$ws.Range("N21").Value = -14936
$ws.Range("H23").Value = 19870.467
$ws.Range("I23").Value = 20773.615
$ws.Range("J23").Value = 14000
$ws.Range("K23").Value = 20773.615
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = -20539.615
$ws.Range("N23").Value = -14468
$ws.Range("H92").Value = 1161.1818
$ws.Range("I92").Value = 1195.5
$ws.Range("J92").Value = 1120
$ws.Range("K92").Value = 1195.5
$ws.Range("L92").Value = 1120
$ws.Range("M92").Value = 52.5
$ws.Range("N92").Value = -3616
$ws.Range("H107").Value = 751.2727
$ws.Range("I107").Value = 719.48
$ws.Range("J107").Value = 850.625
$ws.Range("K107").Value = 719.48
$ws.Range("L107").Value = 850.625
$ws.Range("M107").Value = 1200.52
$ws.Range("N107").Value = -4690.625
$ws.Range("H111").Value = 1918.75
$ws.Range("I111").Value = 1681.4117
$ws.Range("J111").Value = 2495.1428
$ws.Range("K111").Value = 5044.2351
$ws.Range("L111").Value = 7485.428400000001
$ws.Range("M111").Value = -1977.2351
$ws.Range("N111").Value = -13619.4284
$ws.Range("H113").Value = 2002.4546
$ws.Range("I113").Value = 1559.8334
$ws.Range("J113").Value = 2168.4375
$ws.Range("K113").Value = 1559.8334
$ws.Range("L113").Value = 2168.4375
$ws.Range("M113").Value = 1694.1666
$ws.Range("N113").Value = -8676.4375
$ws.Range("H116").Value = 2940.4
$ws.Range("I116").Value = 2915
$ws.Range("K116").Value = 2915
$ws.Range("M116").Value = 527
$ws.Range("H132").Value = 6093.1055
$ws.Range("I132").Value = 6517
$ws.Range("J132").Value = 2490
$ws.Range("K132").Value = 19551
$ws.Range("L132").Value = 7470
$ws.Range("M132").Value = -17021
$ws.Range("N132").Value = -12530
$ws.Range("H135").Value = 734993.2
$ws.Range("I135").Value = 2751.543
$ws.Range("K135").Value = 24763.887
$ws.Range("M135").Value = -22228.887
$ws.Range("H137").Value = 21278038
$ws.Range("I137").Value = 874.4865
$ws.Range("J137").Value = 100003544
$ws.Range("K137").Value = 2623.4595
$ws.Range("L137").Value = 300010632
$ws.Range("M137").Value = -73.45949999999993
$ws.Range("N137").Value = -300015732
$ws.Range("H138").Value = 2322.9333
$ws.Range("I138").Value = 1743.1
$ws.Range("J138").Value = 3482.6
$ws.Range("K138").Value = 5229.299999999999
$ws.Range("L138").Value = 10447.8
$ws.Range("M138").Value = -89.29999999999927
$ws.Range("N138").Value = -20727.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1175.88
$ws.Range("I45").Value = 828.4545000000001
$ws.Range("K45").Value = 828.4545000000001
$ws.Range("M45").Value = -451.4545000000001
$ws.Range("H74").Value = 12824864
$ws.Range("I74").Value = 17857962
$ws.Range("J74").Value = 13341.546
$ws.Range("K74").Value = 17857962
$ws.Range("L74").Value = 13341.546
$ws.Range("M74").Value = -17857088
$ws.Range("N74").Value = -15089.546
$ws.Range("H77").Value = 12824864
$ws.Range("I77").Value = 17857962
$ws.Range("J77").Value = 13341.546
$ws.Range("K77").Value = 89289810
$ws.Range("L77").Value = 66707.73
$ws.Range("M77").Value = -89285442
$ws.Range("N77").Value = -75443.73
$ws.Range("H97").Value = 584.3570999999999
$ws.Range("I97").Value = 550.4286
$ws.Range("J97").Value = 652.2143
$ws.Range("K97").Value = 550.4286
$ws.Range("L97").Value = 652.2143
$ws.Range("M97").Value = -54.42859999999996
$ws.Range("N97").Value = -1644.2143
$ws.Range("H102").Value = 2249.0908
$ws.Range("I102").Value = 2168.5715
$ws.Range("J102").Value = 2390
$ws.Range("K102").Value = 2168.5715
$ws.Range("L102").Value = 2390
$ws.Range("M102").Value = -546.5715
$ws.Range("N102").Value = -5634
$ws.Range("H110").Value = 1510.1111
$ws.Range("J110").Value = 1708
$ws.Range("L110").Value = 1708
$ws.Range("N110").Value = -5798

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1144.1428
$ws.Range("I94").Value = 552.25
$ws.Range("K94").Value = 552.25
$ws.Range("M94").Value = -101.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1537.3096
$ws.Range("I31").Value = 1114.5278
$ws.Range("J31").Value = 4074
$ws.Range("K31").Value = 1114.5278
$ws.Range("L31").Value = 4074
$ws.Range("M31").Value = -819.5278000000001
$ws.Range("N31").Value = -4664
$ws.Range("H34").Value = 1537.3096
$ws.Range("I34").Value = 1114.5278
$ws.Range("J34").Value = 4074
$ws.Range("K34").Value = 1114.5278
$ws.Range("L34").Value = 4074
$ws.Range("M34").Value = -912.5278000000001
$ws.Range("N34").Value = -4478
$ws.Range("H105").Value = 982.75
$ws.Range("I105").Value = 843.3333
$ws.Range("K105").Value = 843.3333
$ws.Range("M105").Value = 903.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 111451.72
$ws.Range("J107").Value = 100356.8
$ws.Range("L107").Value = 301070.4
$ws.Range("N107").Value = -304910.4
$ws.Range("H113").Value = 587.8276
$ws.Range("I113").Value = 615
$ws.Range("J113").Value = 573.5263
$ws.Range("K113").Value = 1845
$ws.Range("L113").Value = 1720.5789
$ws.Range("M113").Value = 325
$ws.Range("N113").Value = -6060.5789
$ws.Range("H114").Value = 1231.25
$ws.Range("I114").Value = 630.2857
$ws.Range("J114").Value = 1832.2142
$ws.Range("K114").Value = 1890.8571
$ws.Range("L114").Value = 5496.642599999999
$ws.Range("M114").Value = 1363.1429
$ws.Range("N114").Value = -12004.6426
$ws.Range("H121").Value = 17857924
$ws.Range("I121").Value = 511.9
$ws.Range("J121").Value = 27778708
$ws.Range("K121").Value = 1535.7
$ws.Range("L121").Value = 83336124
$ws.Range("M121").Value = -225.6999999999998
$ws.Range("N121").Value = -83338744
$ws.Range("H131").Value = 835.03
$ws.Range("I131").Value = 610
$ws.Range("J131").Value = 849.3936
$ws.Range("K131").Value = 1830
$ws.Range("L131").Value = 2548.1808
$ws.Range("M131").Value = 3210
$ws.Range("N131").Value = -12628.1808
$ws.Range("H132").Value = 40001004
$ws.Range("J132").Value = 1735.6666
$ws.Range("L132").Value = 15620.9994
$ws.Range("N132").Value = -20680.9994
$ws.Range("H137").Value = 17847.459
$ws.Range("I137").Value = 2554
$ws.Range("J137").Value = 21872.053
$ws.Range("K137").Value = 7662
$ws.Range("L137").Value = 65616.159
$ws.Range("M137").Value = -2562
$ws.Range("N137").Value = -75816.159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 395.92856
$ws.Range("I107").Value = 371.75
$ws.Range("J107").Value = 428.16666
$ws.Range("K107").Value = 371.75
$ws.Range("L107").Value = 428.16666
$ws.Range("M107").Value = 1548.25
$ws.Range("N107").Value = -4268.16666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1795
$ws.Range("I40").Value = 1595.8334
$ws.Range("J40").Value = 2990
$ws.Range("K40").Value = 1595.8334
$ws.Range("L40").Value = 2990
$ws.Range("M40").Value = -1459.8334
$ws.Range("N40").Value = -3262
$ws.Range("H61").Value = 1349.75
$ws.Range("I61").Value = 1349.2858
$ws.Range("J61").Value = 1350.8334
$ws.Range("K61").Value = 1349.2858
$ws.Range("L61").Value = 1350.8334
$ws.Range("M61").Value = -1147.2858
$ws.Range("N61").Value = -1754.8334
$ws.Range("H93").Value = 1691.2632
$ws.Range("I93").Value = 1280.909
$ws.Range("J93").Value = 2255.5
$ws.Range("K93").Value = 1280.909
$ws.Range("L93").Value = 2255.5
$ws.Range("M93").Value = -32.90900000000011
$ws.Range("N93").Value = -4751.5
$ws.Range("H113").Value = 1349.75
$ws.Range("I113").Value = 1349.2858
$ws.Range("J113").Value = 1350.8334
$ws.Range("K113").Value = 1349.2858
$ws.Range("L113").Value = 1350.8334
$ws.Range("M113").Value = 820.7141999999999
$ws.Range("N113").Value = -5690.8334
$ws.Range("H122").Value = 2650
$ws.Range("I122").Value = 2020
$ws.Range("K122").Value = 6060
$ws.Range("M122").Value = -3610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 13744.134
$ws.Range("I96").Value = 3928.7144
$ws.Range("J96").Value = 22332.625
$ws.Range("K96").Value = 3928.7144
$ws.Range("L96").Value = 22332.625
$ws.Range("M96").Value = -2555.7144
$ws.Range("N96").Value = -25078.625
$ws.Range("H122").Value = 2108.2693
$ws.Range("I122").Value = 1530.75
$ws.Range("K122").Value = 4592.25
$ws.Range("M122").Value = -2142.25
$ws.Range("H126").Value = 2040.7142
$ws.Range("I126").Value = 1823.2142
$ws.Range("J126").Value = 2910.7144
$ws.Range("K126").Value = 5469.642599999999
$ws.Range("L126").Value = 8732.143199999999
$ws.Range("M126").Value = -2999.642599999999
$ws.Range("N126").Value = -13672.1432
